{"js": "// Replace each two-digit multiplication expression in the document with its\n// new value. The mapping is applied as an ordered list of unique\n// old-text -> new-text replacements (each old string occurs exactly once in\n// the document), found via Body.search() and substituted with\n// Range.insertText(..., \"Replace\").\nconst replacements = [\n  [\"14\u00d761=\", \"90\u00d790=\"],\n  [\"68\u00d750=\", \"72\u00d773=\"],\n  [\"78\u00d719=\", \"49\u00d728=\"],\n  [\"19\u00d793=\", \"89\u00d755=\"],\n  [\"48\u00d734=\", \"21\u00d792=\"],\n  [\"32\u00d784=\", \"40\u00d746=\"],\n  [\"70\u00d792=\", \"88\u00d794=\"],\n  [\"29\u00d748=\", \"42\u00d734=\"],\n  [\"97\u00d771=\", \"54\u00d767=\"],\n  [\"91\u00d793=\", \"93\u00d751=\"],\n  [\"15\u00d717=\", \"53\u00d723=\"],\n  [\"55\u00d736=\", \"42\u00d711=\"],\n  [\"46\u00d722=\", \"86\u00d754=\"],\n  [\"12\u00d792=\", \"81\u00d778=\"],\n  [\"83\u00d766=\", \"50\u00d748=\"],\n  [\"27\u00d786=\", \"39\u00d783=\"],\n  [\"69\u00d787=\", \"68\u00d795=\"],\n  [\"50\u00d796=\", \"51\u00d783=\"],\n  [\"18\u00d762=\", \"85\u00d762=\"],\n  [\"87\u00d757=\", \"73\u00d748=\"],\n  [\"31\u00d758=\", \"54\u00d796=\"],\n  [\"87\u00d774=\", \"29\u00d769=\"],\n  [\"27\u00d744=\", \"27\u00d784=\"],\n  [\"73\u00d720=\", \"49\u00d728=\"],\n  [\"93\u00d712=\", \"35\u00d744=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find expression \"${oldText}\" in document body`);\n  }\n\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression in the document with its\n# new value. The mapping is applied as an ordered list of unique\n# old-text -> new-text replacements (each old string occurs exactly once in\n# the document), using Word's Find/Replace engine scoped to the whole story.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"14\u00d761=\", \"90\u00d790=\"),\n  @(\"68\u00d750=\", \"72\u00d773=\"),\n  @(\"78\u00d719=\", \"49\u00d728=\"),\n  @(\"19\u00d793=\", \"89\u00d755=\"),\n  @(\"48\u00d734=\", \"21\u00d792=\"),\n  @(\"32\u00d784=\", \"40\u00d746=\"),\n  @(\"70\u00d792=\", \"88\u00d794=\"),\n  @(\"29\u00d748=\", \"42\u00d734=\"),\n  @(\"97\u00d771=\", \"54\u00d767=\"),\n  @(\"91\u00d793=\", \"93\u00d751=\"),\n  @(\"15\u00d717=\", \"53\u00d723=\"),\n  @(\"55\u00d736=\", \"42\u00d711=\"),\n  @(\"46\u00d722=\", \"86\u00d754=\"),\n  @(\"12\u00d792=\", \"81\u00d778=\"),\n  @(\"83\u00d766=\", \"50\u00d748=\"),\n  @(\"27\u00d786=\", \"39\u00d783=\"),\n  @(\"69\u00d787=\", \"68\u00d795=\"),\n  @(\"50\u00d796=\", \"51\u00d783=\"),\n  @(\"18\u00d762=\", \"85\u00d762=\"),\n  @(\"87\u00d757=\", \"73\u00d748=\"),\n  @(\"31\u00d758=\", \"54\u00d796=\"),\n  @(\"87\u00d774=\", \"29\u00d769=\"),\n  @(\"27\u00d744=\", \"27\u00d784=\"),\n  @(\"73\u00d720=\", \"49\u00d728=\"),\n  @(\"93\u00d712=\", \"35\u00d744=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
